# Remove the blank "spacer" rows that separated trial groups in the data.
# These are blank rows currently at rows 20, 39, 62, 85, 109, 134, 156.
# Deleting them in descending order keeps each row number valid for the
# next deletion (removals below a given row do not shift it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blankRows = @(156, 134, 109, 85, 62, 39, 20)

foreach ($r in $blankRows) {
    $ws.Rows.Item($r).Delete()
}

# Update the view to match the post-edit state.
$ws.Range("A180").Select()

$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 143
